$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: The Nature Conservancy terrestrial ecoregions ---
$ws.Range("A2").Value = "The Nature Conservancy terrestrial ecoregions"
$ws.Range("B2").Value = "Polygon file with all boundaries of terrestrial ecoregions globally"
$ws.Range("C2").Value = "tnc_terr_ecoregions.shp (ESRI shapefile with multiple files)"
$ws.Range("D2").Value = "http://maps.tnc.org/gis_data.html"
$ws.Range("G2").Value = "Olson, D. M. and E. Dinerstein. 2002. The Global 200: Priority ecoregions for global conservation. (PDF file) Annals of the Missouri Botanical Garden 89:125-126."

# --- Row 3: United States administrative boundaries shapefile ---
$ws.Range("A3").Value = "United States administrative boundaries shapefile"
$ws.Range("B3").Value = "Polygon file of the United States' country boundaries based on GADM v2.0"
$ws.Range("C3").Value = "USA_adm0.shp (ESRI shapefile with multiple files)"
$ws.Range("D3").Value = "No longer available. The most recent version is available at https://gadm.org/download_country_v3.html"
$ws.Range("G3").Value = "GADM (Global Administrative Areas Database) [WWW Document], n.d. URL https://gadm.org/ (accessed 4.1.21)."

# --- Row 4: United States county boundaries shapefile ---
$ws.Range("A4").Value = "United States county boundaries shapefile"
$ws.Range("B4").Value = "Polygon file of the United States county boundaries as they existed in 2014"
$ws.Range("C4").Value = "cb_2014_us_county_500k.shp (ESRI shapefile with multiple files)"
$ws.Range("D4").Value = "No longer available. A similar file is available at https://www2.census.gov/geo/tiger/TIGER2014/COUNTY/"
$ws.Range("G4").Value = "U.S. Census Bureau, n.d. TIGER/Line Shapefiles [WWW Document]. The United States Census Bureau. URL https://www.census.gov/geographies/mapping-files/time-series/geo/tiger-line-file.html (accessed 4.1.21)."

# --- Row 5: Global country administrative boundaries shapefile ---
$ws.Range("A5").Value = "Global country administrative boundaries shapefile"
$ws.Range("B5").Value = "Polygon file of all country boundaries as they existed in 2018"
$ws.Range("C5").Value = "ne_50m_admin_0_countries.shp (ESRI shapefile with multiple files)"
$ws.Range("D5").Value = "https://www.naturalearthdata.com/downloads/50m-cultural-vectors/"
$ws.Range("G5").Value = "Natural Earth - Free vector and raster map data at 1:10m, 1:50m, and 1:110m scales, n.d. URL https://www.naturalearthdata.com/ (accessed 4.1.21)."

# --- Row 6: National Land Cover Database 2016, CONUS ---
$ws.Range("A6").Value = "National Land Cover Database 2016, CONUS"
$ws.Range("B6").Value = "Raster at 30m resolution of modeled land cover classes in contiguous United States"
$ws.Range("C6").Value = "NLCD_2016_Land_Cover_L48_20190424.img"
$ws.Range("D6").Value = "https://www.mrlc.gov/data"
$ws.Range("G6").Value = "Dewitz, J., 2019, National Land Cover Database (NLCD) 2016 Products: U.S. Geological Survey data release, https://doi.org/10.5066/P96HHBIE."

# --- Row 7: National Land Cover Database 2016, Alaska ---
$ws.Range("A7").Value = "National Land Cover Database 2016, Alaska"
$ws.Range("B7").Value = "Raster at 30m resolution of modeled land cover classes in Alaska"
$ws.Range("C7").Value = "NLCD_2016_Land_Cover_AK_20200724.img"
$ws.Range("D7").Value = "https://www.mrlc.gov/data/nlcd-2016-land-cover-alaska"
$ws.Range("G7").Value = "Dewitz, J., 2019, National Land Cover Database (NLCD) 2016 Products: U.S. Geological Survey data release, https://doi.org/10.5066/P96HHBIE."

# --- Row 8: NOAA Land Cover Dataset 2001, Hawaii ---
$ws.Range("A8").Value = "NOAA Land Cover Dataset 2001, Hawaii"
$ws.Range("B8").Value = "Raster at 30m resolution of modeled land cover classes in Hawaii"
$ws.Range("C8").Value = "hi_landcover_wimperv_9-30-08_se5.img"
$ws.Range("D8").Value = "https://www.mrlc.gov/data/nlcd-2001-land-cover-hawaii-0"
$ws.Range("G8").Value = "MRLC, 2003. NLCD 2001 Land Cover (HAWAII) | Multi-Resolution Land Characteristics (MRLC) Consortium [WWW Document]. URL https://www.mrlc.gov/data/nlcd-2001-land-cover-hawaii-0 (accessed 4.1.21)."

# --- Row 9: Global pastureland raster layer 2000 ---
$ws.Range("A9").Value = "Global pastureland raster layer 2000"
$ws.Range("B9").Value = "Raster at 1km resolution of global pastureland"
$ws.Range("C9").Value = "pasture.tif"
$ws.Range("D9").Value = "https://sedac.ciesin.columbia.edu/data/set/aglands-pastures-2000/data-download"
$ws.Range("G9").Value = "Ramankutty, N., Evan, A.T., Monfreda, C., Foley, J.A., 2008. Farming the planet: 1. Geographic distribution of global agricultural lands in the year 2000. Global Biogeochemical Cycles 22. https://doi.org/10.1029/2007GB002952"

# --- Row 10: Global cropland dominance raster layer 2010 ---
$ws.Range("A10").Value = "Global cropland dominance raster layer 2010"
$ws.Range("B10").Value = "Raster at 1km resolution of global irrigated and rainfed cropland"
$ws.Range("C10").Value = "GFSAD1KCD.2010.001.2016348142525.tif"
$ws.Range("D10").Value = "https://www.usgs.gov/centers/wgsc/science/global-food-security-support-analysis-data-30-m-gfsad?qt-science_center_objects=4#qt-science_center_objects"
$ws.Range("G10").Value = "Thenkabail, P., Knox, J., Ozdogan, M., Gumma, M., Congalton, R., Wu, Z., Milesi, C., Finkral, A., Marshall, M., Mariotto, I., You, S., Giri, C., Nagler, P. (2016). NASA Making Earth System Data Records for Use in Research Environments (MEaSUREs) Global Food Security Support Analysis Data (GFSAD) Crop Dominance 2010 Global 1 km V001 [Data set]. NASA EOSDIS Land Processes DAAC."

# --- Row 11: United States population raster layer 2010 ---
$ws.Range("A11").Value = "United States population raster layer 2010"
$ws.Range("B11").Value = "Gridded product including population totals from 2010 census at 1 km resolution. Separate files for contiguous USA, Hawaii, Alaska, and Aleutian islands"
$ws.Range("C11").Value = "uspop10.tif, hipop10.tif, akpop10.tif, ehpop10.tif"
$ws.Range("D11").Value = "https://sedac.ciesin.columbia.edu/data/set/usgrid-summary-file1-2010/data-download"
$ws.Range("G11").Value = "Center For International Earth Science Information Network-CIESIN-Columbia University, 2017. U.S. Census Grids (Summary File 1), 2010. https://doi.org/10.7927/H40Z716C"

# --- New hyperlinks for D10 and D11 (D9's hyperlink already exists and is unchanged) ---
$ws.Hyperlinks.Add($ws.Range("D10"), "https://www.usgs.gov/centers/wgsc/science/global-food-security-support-analysis-data-30-m-gfsad", "qt-science_center_objects")
$ws.Hyperlinks.Add($ws.Range("D11"), "https://sedac.ciesin.columbia.edu/data/set/usgrid-summary-file1-2010/data-download")

# --- Column A width widened to fit the newer, longer dataset names ---
$ws.Columns.Item(1).ColumnWidth = 49.7

# --- Active cell moved to J9 ---
$ws.Range("J9").Select()
